# Fixed non-atomic nature of tick calculation.
#
# 1) The fixed "date1" footer text cached on the Slide Master, all 11
#    Custom Layouts, and the Notes Master moves from 9/24/2015 to 10/8/2015
#    (this is the literal text PowerPoint caches for the Header & Footer
#    "Fixed" date field).
# 2) On slide 13 ("Chart Layout"), the four tick-label backing rectangles
#    (Rectangle 4/5/6/7) lose their opaque bg1 fill (now noFill) and three
#    of them are repositioned/resized so the ticks no longer get computed
#    as a separate, non-atomic pass.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer date text: Slide Master + every Custom Layout + Notes Master
# ---------------------------------------------------------------------
$newDate = "10/8/2015"

$master = $p.SlideMaster
$master.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = $newDate

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item("Date Placeholder 2").TextFrame.TextRange.Text = $newDate

# ---------------------------------------------------------------------
# 2) Slide 13 tick rectangles: drop the solid bg1 fill, nudge geometry
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)

# "Rectangle 4" - fill only, geometry untouched
$rect4 = $slide13.Shapes.Item(4)
$rect4.Fill.Visible = $false

# "Rectangle 5" - offset.y moves up slightly, fill removed
$rect5 = $slide13.Shapes.Item(5)
$rect5.Top = 390.0
$rect5.Fill.Visible = $false

# "Rectangle 6" (rotated 270) - widens to the left, fill removed
$rect6 = $slide13.Shapes.Item(6)
$rect6.Left = 98.83016204833984
$rect6.Top = 272.8301696777344
$rect6.Width = 264.3396301269531
$rect6.Fill.Visible = $false

# "Rectangle 7" (rotated 270) - widens to the right, fill removed
$rect7 = $slide13.Shapes.Item(7)
$rect7.Left = 363.0
$rect7.Top = 273.0
$rect7.Width = 264.0
$rect7.Fill.Visible = $false
